$d = $word.ActiveDocument

# "Versi" + "on" -> merge into a single run reading "Version"
# (the Find below only overwrites the "Versi" run's text; "Version" already
# contains the "on" that the following run held, so that following run's
# own text is fixed up next.)
$d.Content.Find.Execute("Versi", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version", 2)

# The run that used to read "on" now reads "onon" when read together with
# the text just inserted before it; trim it back down to "on" so the
# rendered/story text stays "Version" (not "Versionon").
$d.Content.Find.Execute("onon", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "on", 2)

# " 2" -> " 1."  (bump the version number and fold the trailing period that
# used to live in its own run into this one, matching the target markup)
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " 1.", 2)

# The original trailing "." run (the one that sat *after* the _GoBack
# bookmark) is now a redundant leftover period at the very end of the
# story - remove just that one character via a plain Range.Delete() so the
# bookmark in between is left completely untouched.
$end = $d.Content.End
$tail = $d.Range($end - 2, $end - 1)
if ($tail.Text -eq ".") {
    $tail.Delete()
}
